$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.212962962962963
$ws.Range("C2").Value = 0.5123456790123457
$ws.Range("J2").Value = 0.01851851851851852
$ws.Range("P2").Value = 0.1604938271604938
$ws.Range("S2").Value = 0.09567901234567901
$ws.Range("B3").Value = 0.01169590643274854
$ws.Range("C3").Value = 0.02923976608187134
$ws.Range("J3").Value = 0.04678362573099415
$ws.Range("P3").Value = 0.7134502923976608
$ws.Range("S3").Value = 0.1988304093567251
$ws.Range("J4").Value = 0.05357142857142857
$ws.Range("P4").Value = 0.6607142857142857
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.07234042553191489
$ws.Range("D6").Value = 0.01276595744680851
$ws.Range("F6").Value = 0.08085106382978724
$ws.Range("J6").Value = 0.2468085106382979
$ws.Range("O6").Value = 0.00425531914893617
$ws.Range("Q6").Value = 0.1659574468085106
$ws.Range("R6").Value = 0.06808510638297872
$ws.Range("S6").Value = 0.348936170212766
$ws.Range("B7").Value = 0.08280254777070063
$ws.Range("D7").Value = 0.03821656050955414
$ws.Range("E7").Value = 0.006369426751592357
$ws.Range("F7").Value = 0.05095541401273886
$ws.Range("J7").Value = 0.1146496815286624
$ws.Range("O7").Value = 0.02547770700636943
$ws.Range("Q7").Value = 0.1910828025477707
$ws.Range("R7").Value = 0.08280254777070063
$ws.Range("S7").Value = 0.4076433121019108
$ws.Range("B8").Value = 0.0996309963099631
$ws.Range("D8").Value = 0.02767527675276753
$ws.Range("F8").Value = 0.04797047970479705
$ws.Range("J8").Value = 0.1365313653136531
$ws.Range("O8").Value = 0.01845018450184502
$ws.Range("Q8").Value = 0.1752767527675277
$ws.Range("R8").Value = 0.0940959409594096
$ws.Range("S8").Value = 0.4003690036900369
$ws.Range("B9").Value = 0.1077844311377246
$ws.Range("D9").Value = 0.03592814371257485
$ws.Range("F9").Value = 0.03592814371257485
$ws.Range("J9").Value = 0.08383233532934131
$ws.Range("O9").Value = 0.01197604790419162
$ws.Range("Q9").Value = 0.1976047904191617
$ws.Range("R9").Value = 0.08383233532934131
$ws.Range("S9").Value = 0.4431137724550898
$ws.Range("B10").Value = 0.1004126547455296
$ws.Range("D10").Value = 0.01788170563961486
$ws.Range("F10").Value = 0.07634112792297111
$ws.Range("J10").Value = 0.1272352132049518
$ws.Range("O10").Value = 0.01306740027510316
$ws.Range("Q10").Value = 0.2290233837689133
$ws.Range("R10").Value = 0.0859697386519945
$ws.Range("S10").Value = 0.3500687757909216
$ws.Range("G11").Value = 0.1406844106463878
$ws.Range("J11").Value = 0.1178707224334601
$ws.Range("K11").Value = 0.1977186311787072
$ws.Range("L11").Value = 0.5209125475285171
$ws.Range("S11").Value = 0.02281368821292776
$ws.Range("G12").Value = 0.6643356643356644
$ws.Range("J12").Value = 0.2097902097902098
$ws.Range("K12").Value = 0.03496503496503497
$ws.Range("L12").Value = 0.04895104895104895
$ws.Range("S12").Value = 0.04195804195804196
$ws.Range("G13").Value = 0.5952380952380952
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.01158301158301158
$ws.Range("H15").Value = 0.2277992277992278
$ws.Range("I15").Value = 0.05405405405405406
$ws.Range("J15").Value = 0.3822393822393823
$ws.Range("K15").Value = 0.02702702702702703
$ws.Range("M15").Value = 0.01544401544401544
$ws.Range("O15").Value = 0.09652509652509653
$ws.Range("S15").Value = 0.1853281853281853
$ws.Range("F16").Value = 0.01485148514851485
$ws.Range("H16").Value = 0.1683168316831683
$ws.Range("I16").Value = 0.0891089108910891
$ws.Range("J16").Value = 0.3811881188118812
$ws.Range("K16").Value = 0.07425742574257425
$ws.Range("M16").Value = 0.0198019801980198
$ws.Range("O16").Value = 0.06930693069306931
$ws.Range("S16").Value = 0.1831683168316832
$ws.Range("F17").Value = 0.02466793168880456
$ws.Range("H17").Value = 0.206831119544592
$ws.Range("I17").Value = 0.07400379506641366
$ws.Range("J17").Value = 0.4326375711574952
$ws.Range("K17").Value = 0.07210626185958255
$ws.Range("M17").Value = 0.01707779886148008
$ws.Range("N17").Value = 0.00189753320683112
$ws.Range("O17").Value = 0.06261859582542695
$ws.Range("S17").Value = 0.1081593927893738
$ws.Range("F18").Value = 0.01818181818181818
$ws.Range("H18").Value = 0.1909090909090909
$ws.Range("I18").Value = 0.08181818181818182
$ws.Range("J18").Value = 0.45
$ws.Range("K18").Value = 0.07727272727272727
$ws.Range("M18").Value = 0.02272727272727273
$ws.Range("O18").Value = 0.05909090909090909
$ws.Range("S18").Value = 0.1
$ws.Range("F19").Value = 0.008869179600886918
$ws.Range("H19").Value = 0.2254249815225425
$ws.Range("I19").Value = 0.05838876570583888
$ws.Range("J19").Value = 0.3909830007390983
$ws.Range("K19").Value = 0.09312638580931264
$ws.Range("M19").Value = 0.01699926090169993
$ws.Range("N19").Value = 0.002217294900221729
$ws.Range("O19").Value = 0.07908351810790835
$ws.Range("S19").Value = 0.1249076127124908
